$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new data row (row 4) mirroring the structure of existing rows
$ws.Cells.Item(4, 1).Value = 42605.886620370373
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(4, 2).Value = -28
$ws.Cells.Item(4, 3).Value = 49
$ws.Cells.Item(4, 4).Value = 48
$ws.Cells.Item(4, 5).Value = 22
$ws.Cells.Item(4, 6).Value = 77
$ws.Cells.Item(4, 7).Value = 15700
$ws.Cells.Item(4, 8).Value = 8713
$ws.Cells.Item(4, 9).Value = 1089
$ws.Cells.Item(4, 10).Value = 105
$ws.Cells.Item(4, 11).Value = 104
$ws.Cells.Item(4, 12).Value = 6
$ws.Cells.Item(4, 13).Value = 21
$ws.Cells.Item(4, 14).Value = "Bag"
